# Editor with cell number
# Applies the "difficulty" worksheet rework:
#  - New Level numbers (column B) for rows 4-37
#  - New Rows/Cols/Nb Shapes data (columns C, D, F) for rows 4-13
#  - Clears the old Rows/Cols/Nb Shapes/Difficulty data for rows 14-18
#  - Removes the old "Difficulty" text column (G) entirely
#  - Adds a small legend/side table in columns K:M (rows 5-11)
#  - Updates the active cell selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Column B: sequential Level numbers for rows 4-37 ----
for ($r = 4; $r -le 37; $r++) {
    $ws.Cells.Item($r, 2).Value = $r - 3
}

# ---- Columns C (Rows), D (Cols), F (Nb Shapes) for rows 4-13 ----
$newData = @{
    4  = @(5, 5, 3)
    5  = @(5, 5, 4)
    6  = @(5, 5, 5)
    7  = @(6, 6, 6)
    8  = @(6, 6, 7)
    9  = @(6, 6, 8)
    10 = @(7, 7, 6)
    11 = @(7, 7, 7)
    12 = @(7, 7, 8)
    13 = @(7, 7, 9)
}
foreach ($r in $newData.Keys) {
    $vals = $newData[$r]
    $ws.Cells.Item($r, 3).Value = $vals[0]
    $ws.Cells.Item($r, 4).Value = $vals[1]
    $ws.Cells.Item($r, 6).Value = $vals[2]
}

# ---- Clear old Rows/Cols/Nb Shapes data left over in rows 14-18 ----
$ws.Range("C14:D18").ClearContents() | Out-Null
$ws.Range("F14:F18").ClearContents() | Out-Null

# ---- Remove the old "Difficulty" text column entirely (rows 4-13) ----
$ws.Range("G4:G13").ClearContents() | Out-Null

# ---- New legend/side table in columns K:L:M ----
$ws.Range("L5").Value = "12 shapes"
$ws.Range("K7").Value = "Normal"
$ws.Range("L6").Value = "shapes"
$ws.Range("M6").Value = "nbcel"
$ws.Range("M7").Value = "50/70"
$ws.Range("K8").Value = "Easy"
$ws.Range("K9").Value = "Meidum"
$ws.Range("K10").Value = "Hard"
$ws.Range("K11").Value = "Master"
$ws.Range("L7").Value = 8
$ws.Range("L8").Value = 6
$ws.Range("L9").Value = 9
$ws.Range("L10").Value = 11
$ws.Range("L11").Value = 12
$ws.Range("M10").Value = 80

# ---- Update selection to match the author's last-active cell ----
$ws.Range("L14").Select() | Out-Null
